# fix: solucion o cambio
# Corrige las cantidades realizadas de las actividades asignadas de enero
# 2026 (filas 9 y 10) y registra el medio de verificacion de la fila 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Participacion en el desarrollo y revision del AST": solo se realizo 1 de 3
$ws.Range("F9").Value = 1

# "Difusion de charla Diaria": solo se realizo 1 de 2
$ws.Range("F10").Value = 1

# Medio de verificacion de "Otros:" (fila 18) - se ingresa como texto "1"
$ws.Range("G18").Value = "1"

# Deja la seleccion activa sobre la ultima celda editada
$ws.Range("F10").Select()
